$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1818.7778
$ws.Range("I28").Value = 1122.8
$ws.Range("J28").Value = 2086.4614
$ws.Range("K28").Value = 1122.8
$ws.Range("L28").Value = 2086.4614
$ws.Range("M28").Value = -637.8
$ws.Range("N28").Value = -3056.4614
$ws.Range("H38").Value = 4808
$ws.Range("J38").Value = 13055.5
$ws.Range("L38").Value = 39166.5
$ws.Range("N38").Value = -39910.5
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").ClearContents()
$ws.Range("N48").Value = 0
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").ClearContents()
$ws.Range("N56").Value = 0
$ws.Range("H97").Value = 1937.6666
$ws.Range("J97").Value = 2125.2
$ws.Range("L97").Value = 6375.599999999999
$ws.Range("N97").Value = -7367.599999999999
$ws.Range("H112").Value = 1417.6041
$ws.Range("J112").Value = 1425.1111
$ws.Range("L112").Value = 4275.3333
$ws.Range("N112").Value = -6491.3333
$ws.Range("H125").Value = 2849.9333
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 2910.6428
$ws.Range("K125").Value = 18000
$ws.Range("L125").Value = 26195.7852
$ws.Range("M125").Value = -15540
$ws.Range("N125").Value = -31115.7852
$ws.Range("H132").Value = 15157.393
$ws.Range("I132").Value = 7135.5386
$ws.Range("J132").Value = 22109.666
$ws.Range("K132").Value = 21406.6158
$ws.Range("L132").Value = 66328.99800000001
$ws.Range("M132").Value = -18876.6158
$ws.Range("N132").Value = -71388.99800000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 5250
$ws.Range("I30").Value = 500
$ws.Range("K30").Value = 500
$ws.Range("M30").Value = -350
$ws.Range("H32").Value = 13504.58
$ws.Range("J32").Value = 10480.375
$ws.Range("L32").Value = 10480.375
$ws.Range("N32").Value = -11054.375
$ws.Range("H36").Value = 4801.1113
$ws.Range("I36").Value = 4651.25
$ws.Range("J36").Value = 6000
$ws.Range("K36").Value = 4651.25
$ws.Range("L36").Value = 6000
$ws.Range("M36").Value = -4305.25
$ws.Range("N36").Value = -6692
$ws.Range("H61").Value = 6823.8535
$ws.Range("I61").Value = 8042.0645
$ws.Range("K61").Value = 8042.0645
$ws.Range("M61").Value = -7830.0645
$ws.Range("H97").Value = 407.86957
$ws.Range("I97").Value = 411.5238
$ws.Range("K97").Value = 411.5238
$ws.Range("M97").Value = 84.47620000000001
$ws.Range("H132").Value = 5067.881
$ws.Range("I132").Value = 2427.65
$ws.Range("J132").Value = 7468.091
$ws.Range("K132").Value = 7282.950000000001
$ws.Range("L132").Value = 22404.273
$ws.Range("M132").Value = -4752.950000000001
$ws.Range("N132").Value = -27464.273
$ws.Range("H136").Value = 6823.8535
$ws.Range("I136").Value = 8042.0645
$ws.Range("K136").Value = 24126.1935
$ws.Range("M136").Value = -21576.1935

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1489504.4
$ws.Range("I99").Value = 2315940.2
$ws.Range("K99").Value = 2315940.2
$ws.Range("M99").Value = -2314442.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 15890408
$ws.Range("I132").Value = 18532732
$ws.Range("J132").Value = 36466.668
$ws.Range("K132").Value = 55598196
$ws.Range("L132").Value = 109400.004
$ws.Range("M132").Value = -55595666
$ws.Range("N132").Value = -114460.004
$ws.Range("H134").Value = 1680.4286
$ws.Range("I134").Value = 1350.8379
$ws.Range("K134").Value = 4052.5137
$ws.Range("M134").Value = -1517.5137

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 19581
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 19581
$ws.Range("K64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("M64").Value = 58743
$ws.Range("N64").Value = -59283
$ws.Range("H67").Value = 19581
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 19581
$ws.Range("K67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("M67").Value = 58743
$ws.Range("N67").Value = -60615
$ws.Range("H113").Value = 823.6667
$ws.Range("J113").Value = 960.5
$ws.Range("L113").Value = 2881.5
$ws.Range("N113").Value = -7221.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1933.6666
$ws.Range("I113").Value = 1801
$ws.Range("K113").Value = 1801
$ws.Range("M113").Value = 369
$ws.Range("H122").Value = 412128.38
$ws.Range("I122").Value = 689780
$ws.Range("K122").Value = 2069340
$ws.Range("M122").Value = -2066890
$ws.Range("H123").Value = 48979.25
$ws.Range("J123").Value = 48979.25
$ws.Range("L123").Value = 48979.25
$ws.Range("N123").Value = -53879.25
$ws.Range("H127").Value = 80233.336
$ws.Range("J127").Value = 80233.336
$ws.Range("L127").Value = 80233.336
$ws.Range("N127").Value = -90153.336
$ws.Range("H132").Value = 141037.14
$ws.Range("I132").Value = 256945.12
$ws.Range("J132").Value = 8570.857
$ws.Range("K132").Value = 770835.36
$ws.Range("L132").Value = 25712.571
$ws.Range("M132").Value = -768305.36
$ws.Range("N132").Value = -30772.571

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2824.524
$ws.Range("I93").Value = 3129.889
$ws.Range("J93").Value = 992.3333
$ws.Range("K93").Value = 3129.889
$ws.Range("L93").Value = 992.3333
$ws.Range("M93").Value = -1881.889
$ws.Range("N93").Value = -3488.3333
$ws.Range("H132").Value = 4360.8
$ws.Range("I132").Value = 4360.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13082.4
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -10552.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5279.364
$ws.Range("J96").Value = 5328.4287
$ws.Range("L96").Value = 5328.4287
$ws.Range("N96").Value = -8074.4287
$ws.Range("H126").Value = 4597.5713
$ws.Range("I126").Value = 4698.8
$ws.Range("J126").Value = 4541.3335
$ws.Range("K126").Value = 14096.4
$ws.Range("L126").Value = 13624.0005
$ws.Range("M126").Value = -11626.4
$ws.Range("N126").Value = -18564.0005
$ws.Range("H132").Value = 11908345
$ws.Range("I132").Value = 1139.9412
$ws.Range("K132").Value = 3419.8236
$ws.Range("M132").Value = -889.8235999999997
$ws.Range("H136").Value = 8327.35
$ws.Range("I136").Value = 2647.889
$ws.Range("K136").Value = 7943.667
$ws.Range("M136").Value = -5393.667
